$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B2").Value = "Espace pour revenir au début des instructions`nBienvenue à la course des champions `nMerci d'avoir rejoint le jury !"
$ws.Range("B2").Select()
